# [AFG] added final excel sheets for Afghanistan
#
# 1) Clear the stray empty inline-string cells in column B of "ODI Batting"
#    (rows 2,5,6,7,8,10,11,12 had an empty <c t="inlineStr"/> cell which should
#    simply not exist).
# 2) Add two brand-new sheets, "ODI Batting Extra" and "ODI Bowling Extra",
#    placed right after "ODI Bowling", with their header + data rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: force a value to be written as TEXT (so numeric-looking strings like
# "3709" or "0.90%" stay strings instead of turning into real numbers), then
# reset the style back to Normal so we don't leave a stray "@" text format on
# the cell.
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------------
# 1) ODI Batting: drop the empty column-B placeholder cells.
# ---------------------------------------------------------------------------
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$emptyBRows = @(2, 5, 6, 7, 8, 10, 11, 12)
foreach ($r in $emptyBRows) {
    $odiBatting.Cells.Item($r, 2).Value = $null
}

# ---------------------------------------------------------------------------
# 2) Add "ODI Batting Extra" after "ODI Bowling"
# ---------------------------------------------------------------------------
$odiBowling = $wb.Worksheets.Item("ODI Bowling")
$battingExtra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $odiBowling)
$battingExtra.Name = "ODI Batting Extra"

$battingHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 0; $c -lt $battingHeaders.Length; $c++) {
    $battingExtra.Cells.Item(1, $c + 1).Value = $battingHeaders[$c]
}
# copy the bold/bordered header look used by the other sheets
$odiBatting.Range("A1:F1").Copy()
$battingExtra.Range("A1:F1").PasteSpecial(-4122) | Out-Null

# row -> @{ col letter = value }; numbers are real numbers, everything else is text
$battingRows = @{
    2  = @{ A = "3709"; F = "NO" }
    3  = @{ A = "3994"; B = 11;  C = "0"; D = "0"; E = "0.90%"; F = "NO" }
    4  = @{ A = "3996"; B = 9;   F = "NO" }
    5  = @{ A = "4010"; F = "NO" }
    6  = @{ A = "4046"; B = 10;  F = "NO" }
    7  = @{ A = "4528"; B = 10;  F = "NO" }
    8  = @{ A = "4530"; F = "NO" }
    9  = @{ A = "4538"; B = 11;  F = "NO" }
    10 = @{ A = "4582" }
    11 = @{ A = "4585" }
    12 = @{ A = "4588" }
}
$numericCols = @("B")

foreach ($r in $battingRows.Keys) {
    $rowData = $battingRows[$r]
    foreach ($col in $rowData.Keys) {
        $cell = $battingExtra.Range("$col$r")
        $val = $rowData[$col]
        if ($numericCols -contains $col) {
            $cell.Value = $val
        } else {
            Set-TextValue $cell $val
        }
    }
}

# ---------------------------------------------------------------------------
# 3) Add "ODI Bowling Extra" after "ODI Batting Extra"
# ---------------------------------------------------------------------------
$bowlingExtra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $battingExtra)
$bowlingExtra.Name = "ODI Bowling Extra"

$bowlingHeaders = @("MATCH_CODE", "MAIDEN_OVERS", "PERCENT_WICKETS_OF_ALL")
for ($c = 0; $c -lt $bowlingHeaders.Length; $c++) {
    $bowlingExtra.Cells.Item(1, $c + 1).Value = $bowlingHeaders[$c]
}
$odiBatting.Range("A1:C1").Copy()
$bowlingExtra.Range("A1:C1").PasteSpecial(-4122) | Out-Null

$bowlingRows = @{
    2  = @{ A = "3709"; B = "";  C = "" }
    3  = @{ A = "3994"; B = "0"; C = "20.00%" }
    4  = @{ A = "3996"; B = "0"; C = "" }
    5  = @{ A = "4010"; B = "0"; C = "30.00%" }
    6  = @{ A = "4528"; B = "1"; C = "20.00%" }
    7  = @{ A = "4530"; B = "0"; C = "10.00%" }
    8  = @{ A = "4538"; B = "";  C = "" }
    9  = @{ A = "4582"; B = "2"; C = "10.00%" }
    10 = @{ A = "4585"; B = "0"; C = "30.00%" }
    11 = @{ A = "4588"; B = "1"; C = "10.00%" }
}

foreach ($r in $bowlingRows.Keys) {
    $rowData = $bowlingRows[$r]
    foreach ($col in $rowData.Keys) {
        $cell = $bowlingExtra.Range("$col$r")
        $val = $rowData[$col]
        if ($val -eq "") {
            # leave the cell blank - nothing to write
            continue
        }
        Set-TextValue $cell $val
    }
}
